$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, shifting existing rows 15-19 down to 16-20
$ws.Rows.Item(15).Insert()

# Populate the new row 15 with fresh data
$ws.Cells.Item(15, 1).Value = 8
$ws.Cells.Item(15, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(15, 3).Value = "Coquimbo"
$ws.Cells.Item(15, 4).Value = 45180
$ws.Cells.Item(15, 5).Value = 4
$ws.Cells.Item(15, 6).Value = 100114002
$ws.Cells.Item(15, 7).Value = "Camote"
$ws.Cells.Item(15, 8).Value = "Sin especificar"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 400
$ws.Cells.Item(15, 11).Value = 16500
$ws.Cells.Item(15, 12).Value = 17000
$ws.Cells.Item(15, 13).Value = 16750
$ws.Cells.Item(15, 14).Value = "$/malla 18 kilos"
$ws.Cells.Item(15, 15).Value = "Perú"
$ws.Cells.Item(15, 16).Value = 931
$ws.Cells.Item(15, 17).Value = 18
$ws.Cells.Item(15, 18).Value = "Hortaliza"
